$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H (copy formatting from the neighboring header cell)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Save flag values for rows 2-11 (1 if it was a save, else 0)
$saveValues = @(0, 0, 1, 0, 0, 0, 0, 0, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
